# update constraint & bounds logic
#
# The "Upper" bound (column C) for four asset classes - 15+ STRIPS (row 3),
# Long Corporate (row 4), Equity (row 6) and Liquid Alternatives (row 7) -
# is tightened from 1.02 to 1 on every sheet (IBT, Pension, Retirement).
# The saved cell selection on each sheet is also updated to reflect the
# new state of the workbook.

$wb = $excel.ActiveWorkbook

$sheetNames = @("IBT", "Pension", "Retirement")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("C3").Value = 1
    $ws.Range("C4").Value = 1
    $ws.Range("C6").Value = 1
    $ws.Range("C7").Value = 1
}

# Update the saved selection on each sheet.
$wsIBT = $wb.Worksheets.Item("IBT")
$wsPension = $wb.Worksheets.Item("Pension")
$wsRetirement = $wb.Worksheets.Item("Retirement")

$wsPension.Range("C2").Select() | Out-Null
$wsRetirement.Range("C2").Select() | Out-Null

# IBT stays the active/selected tab, with C1 as its selection.
$wsIBT.Activate() | Out-Null
$wsIBT.Range("C1").Select() | Out-Null
